# LAHomesForRent sheet (the active sheet) previously listed a single
# "Los Angeles" row. The commit adds two more rows (Nashville, Phoenix)
# and tidies up row 2 (count 6 -> 5, and the LA text loses its stray
# trailing space).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing formatting of row 2 (A2 -> style "2" text count column,
# B2 -> style "3" text/Menlo city column) down into the two new rows so the
# new cells pick up the same number format / font instead of Excel
# auto-detecting them as numbers.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

# Update row 2 and fill in the two new rows.
$ws.Range("A2").Value = "5"
$ws.Range("B2").Value = "Homes for Rent in Los Angeles, CA"

$ws.Range("A3").Value = "6"
$ws.Range("B3").Value = "Homes for Rent in Nashville, TN"

$ws.Range("A4").Value = "9"
$ws.Range("B4").Value = "Homes for Rent in Phoenix, AZ"

# Match the saved selection left behind in the sheet.
$ws.Range("B9").Select() | Out-Null
